$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("F13").Value = 14
$ws.Range("F19").Value = 27

# Add new row 24 with formula referencing the existing sum in F23
$ws.Range("F24").Formula = "=F23-200"

# Update the selection to match the final state
$ws.Range("F14").Select()
